$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.289535760879517
$ws.Range("B1").Value = 2.606353998184204
$ws.Range("C1").Value = 2.063186168670654
$ws.Range("D1").Value = 1.943684458732605
$ws.Range("E1").Value = 1.71659529209137
